# Change D4:D22 on the "Test Cases" sheet from "Y" to "N"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

for ($r = 4; $r -le 22; $r++) {
    $ws.Cells.Item($r, 4).Value = "N"
}
